$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the now-stale WH 1202097 / PO 1717213 line) — rows below shift up
$ws.Rows.Item(2).Delete()

# Restore the active-cell selection to C5 to match the saved view state
$ws.Range("C5").Select()
